$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate column O's formatting (styles, number formats, etc.) into a
# brand-new column P by copying the whole column and inserting it, then
# overwrite the copied values with the real 2022-column data.
$ws.Columns("O").Copy()
$ws.Columns("P").Insert(-4161)
$excel.CutCopyMode = 0

$ws.Range("P4").Value = 2022
$ws.Range("P5").Value = 1
$ws.Range("P6").Value = "-"
$ws.Range("P7").Value = "-"
$ws.Range("P8").Value = "-"
$ws.Range("P9").Value = "-"
$ws.Range("P10").Value = "-"
$ws.Range("P11").Value = "-"
$ws.Range("P12").Value = 1
$ws.Range("P13").Value = "-"
$ws.Range("P14").Value = "-"

# Move the active selection to match the author's last selection.
$ws.Range("O21:O22").Select()
